$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.425"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05890"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.442"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.542"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8110"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9535"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07451"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03264"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03063"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09335"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.856"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001575"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04684"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005958"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.005865"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.004895"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00006809"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.592"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").Value = "'0.0002287"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03931"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006191"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Value = "'0.009258"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005210"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.7310"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.002364"
$ws.Range("D48").Style = "Normal"
